$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E2").ClearContents()

$ws.Range("E2").Select()
